$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = '&__Invariant when using asp-for=""'
$ws.Range("A11").Value = "why pass data using ViewBag instead of passing model?"

$ws.Range("A11").Select()
